# Update "Forecast Comparison" sheet: insert a Week_Start_Date column,
# shorten the Week labels, correct two MyForecast values, and store
# is_holiday_week as a real boolean. Then refresh the dependent total
# on the Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Insert the new "Week_Start_Date" column before the ASIN column (B) ---
$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "Week_Start_Date"

# Store the dates as plain text (e.g. "2025-01-05"), not Excel date serials
$ws.Range("B2:B17").NumberFormat = "@"

# --- Per-week data: Week label, Week start date, MyForecast, Amazon Mean,
#     Amazon P70, Amazon P80, Amazon P90 ---
$weekRows = @(
    @("W1",  "2025-01-05", 17, 10, 10, 16, 28),
    @("W2",  "2025-01-12", 31, 11, 10, 18, 32),
    @("W3",  "2025-01-19", 36, 13, 13, 22, 39),
    @("W4",  "2025-01-26", 33, 14, 14, 24, 41),
    @("W5",  "2025-02-02", 35, 14, 14, 23, 41),
    @("W6",  "2025-02-09", 37, 14, 13, 22, 40),
    @("W7",  "2025-02-16", 34, 14, 13, 23, 42),
    @("W8",  "2025-02-23", 38, 14, 11, 21, 41),
    @("W9",  "2025-03-02", 36, 13, 12, 21, 39),
    @("W10", "2025-03-09", 34, 13, 10, 20, 39),
    @("W11", "2025-03-16", 32, 13, 11, 20, 39),
    @("W12", "2025-03-23", 33, 13, 11, 20, 40),
    @("W13", "2025-03-30", 33, 12, 10, 19, 37),
    @("W14", "2025-04-06", 32, 12,  9, 18, 36),
    @("W15", "2025-04-13", 31, 12, 10, 19, 38),
    @("W16", "2025-04-20", 31, 12,  9, 18, 36)
)

$r = 2
foreach ($row in $weekRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]               # A: Week (no leading zero)
    $ws.Cells.Item($r, 2).Value = $row[1]                # B: Week_Start_Date
    $ws.Cells.Item($r, 3).Value = "B0BZ17BQ4Z"           # C: ASIN
    $ws.Cells.Item($r, 4).Value = $row[2]                # D: MyForecast
    $ws.Cells.Item($r, 5).Value = $row[3]                # E: Amazon Mean Forecast
    $ws.Cells.Item($r, 6).Value = $row[4]                # F: Amazon P70 Forecast
    $ws.Cells.Item($r, 7).Value = $row[5]                # G: Amazon P80 Forecast
    $ws.Cells.Item($r, 8).Value = $row[6]                # H: Amazon P90 Forecast
    $ws.Cells.Item($r, 9).Value = "B650M K"              # I: Product Title
    $ws.Cells.Item($r, 10).Value = $false                # J: is_holiday_week (boolean)
    $r = $r + 1
}

# --- Refresh the dependent "Total Forecast (8 Weeks)" figure on Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B10").NumberFormat = "@"
$summary.Range("B10").Value = "261"
